$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Experimental flag: "false" -> "true" (keep as literal text, not a Boolean)
$ws.Range("B7").Formula = '="true"'
$ws.Range("B7").Copy()
$ws.Range("B7").PasteSpecial(-4163)

# Date updated
$ws.Range("B8").Value = "2023-02-16T14:43:10-06:00"

# Case Sensitive value now set to "true" (was blank) - keep as literal text
$ws.Range("B14").Formula = '="true"'
$ws.Range("B14").Copy()
$ws.Range("B14").PasteSpecial(-4163)

$excel.CutCopyMode = 0
